$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Label" header column (H), matching the header style used by the
# other header cells (bold, centered, thin border) by copying the format
# from the neighboring header cell (G1) without disturbing the new value.
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)   # xlPasteFormats

# Build map of row -> label value based on patient name in column A
$labelMap = @{
    "Control 39" = 0
    "Control 17" = 0
    "Control 23" = 0
    "Control 27" = 0
    "Control 8"  = 0
    "MDD 36" = 1
    "MDD 10" = 1
    "MDD 39" = 1
    "MDD 14" = 1
    "MDD 18" = 1
}

for ($r = 2; $r -le 21; $r++) {
    $name = $ws.Cells.Item($r, 1).Value2
    if ($labelMap.ContainsKey($name)) {
        $ws.Cells.Item($r, 8).Value = $labelMap[$name]
    }
}
